# Sheets scheduled-runner update
# Refreshes cached market-board pricing columns (H:N) for a batch of Leve
# rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
#   H = currentAveragePrice        K = LevePriceNQ
#   I = currentAveragePriceNQ      L = LevePriceHQ
#   J = currentAveragePriceHQ      M = LeveProfitNQ
#                                   N = LeveProfitHQ

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 29: Weak Blinding Potion
$ws.Cells.Item(29, 8).Value = 2333.3333
$ws.Cells.Item(29, 10).Value = 2333.3333
$ws.Cells.Item(29, 12).Value = 6999.999899999999
$ws.Cells.Item(29, 14).Value = -7561.999899999999
# Row 31: Weak Silencing Potion
$ws.Cells.Item(31, 8).Value = 1000
$ws.Cells.Item(31, 9).Value = 2000
$ws.Cells.Item(31, 10).Value = 500
$ws.Cells.Item(31, 11).Value = 6000
$ws.Cells.Item(31, 12).Value = 1500
$ws.Cells.Item(31, 13).Value = -5770
$ws.Cells.Item(31, 14).Value = -1960
# Row 58: Mega-Potion of Vitality
$ws.Cells.Item(58, 8).Value = 281.5
$ws.Cells.Item(58, 10).Value = 96
$ws.Cells.Item(58, 12).Value = 288
$ws.Cells.Item(58, 14).Value = -588
# Row 74: Wing Glue
$ws.Cells.Item(74, 8).Value = 5399.8
$ws.Cells.Item(74, 9).Value = 5399.8
$ws.Cells.Item(74, 11).Value = 5399.8
$ws.Cells.Item(74, 13).Value = -4463.8
# Row 76: Enchanted Hardsilver Ink
$ws.Cells.Item(76, 8).Value = 45457644
$ws.Cells.Item(76, 9).Value = 100002744
$ws.Cells.Item(76, 11).Value = 100002744
$ws.Cells.Item(76, 13).Value = -100002429
# Row 77: Wing Glue
$ws.Cells.Item(77, 8).Value = 5399.8
$ws.Cells.Item(77, 9).Value = 5399.8
$ws.Cells.Item(77, 11).Value = 26999
$ws.Cells.Item(77, 13).Value = -22319
# Row 79: Enchanted Hardsilver Ink
$ws.Cells.Item(79, 8).Value = 45457644
$ws.Cells.Item(79, 9).Value = 100002744
$ws.Cells.Item(79, 11).Value = 100002744
$ws.Cells.Item(79, 13).Value = -100001652
# Row 112: Superior Spiritbond Potion
$ws.Cells.Item(112, 8).Value = 1273.2916
$ws.Cells.Item(112, 10).Value = 1447.3684
$ws.Cells.Item(112, 12).Value = 4342.1052
$ws.Cells.Item(112, 14).Value = -6558.1052
# Row 115: Competent Craftsman's Syrup
$ws.Cells.Item(115, 8).Value = 1342.5
$ws.Cells.Item(115, 9).Value = 1342.5
$ws.Cells.Item(115, 11).Value = 4027.5
$ws.Cells.Item(115, 13).Value = -2460.5
# Row 125: Grade 5 Dexterity Alkahest
$ws.Cells.Item(125, 8).Value = 9097.6
$ws.Cells.Item(125, 10).Value = 4750
$ws.Cells.Item(125, 12).Value = 42750
$ws.Cells.Item(125, 14).Value = -47670

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Steel Ingot
$ws.Cells.Item(32, 8).Value = 4037.442
$ws.Cells.Item(32, 9).Value = 1598.8551
$ws.Cells.Item(32, 10).Value = 13935.235
$ws.Cells.Item(32, 11).Value = 1598.8551
$ws.Cells.Item(32, 12).Value = 13935.235
$ws.Cells.Item(32, 13).Value = -1311.8551
$ws.Cells.Item(32, 14).Value = -14509.235

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Adamantite Nugget
$ws.Cells.Item(86, 8).Value = 3741.2083
$ws.Cells.Item(86, 9).Value = 2838.111
$ws.Cells.Item(86, 10).Value = 6450.5
$ws.Cells.Item(86, 11).Value = 2838.111
$ws.Cells.Item(86, 12).Value = 6450.5
$ws.Cells.Item(86, 13).Value = -1715.111
$ws.Cells.Item(86, 14).Value = -8696.5
# Row 88: Adamantite Zweihander
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 13).ClearContents()  # M88
$ws.Cells.Item(88, 14).ClearContents()  # N88
# Row 89: Adamantite Nugget
$ws.Cells.Item(89, 8).Value = 3741.2083
$ws.Cells.Item(89, 9).Value = 2838.111
$ws.Cells.Item(89, 10).Value = 6450.5
$ws.Cells.Item(89, 11).Value = 14190.555
$ws.Cells.Item(89, 12).Value = 32252.5
$ws.Cells.Item(89, 13).Value = -8574.555
$ws.Cells.Item(89, 14).Value = -43484.5
# Row 91: Adamantite Zweihander
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 13).ClearContents()  # M91
$ws.Cells.Item(91, 14).ClearContents()  # N91
# Row 107: Deepgold Nugget
$ws.Cells.Item(107, 8).Value = 3615.75
$ws.Cells.Item(107, 10).Value = 5555
$ws.Cells.Item(107, 12).Value = 5555
$ws.Cells.Item(107, 14).Value = -9395
# Row 132: Mountain Chromite Twinfangs
$ws.Cells.Item(132, 8).Value = 33166.297
$ws.Cells.Item(132, 10).Value = 33166.297
$ws.Cells.Item(132, 12).Value = 33166.297
$ws.Cells.Item(132, 14).Value = -43286.297

$ws = $wb.Worksheets.Item("CRP")
# Row 99: Pine Lumber
$ws.Cells.Item(99, 8).Value = 7911060.5
$ws.Cells.Item(99, 9).Value = 11112989
$ws.Cells.Item(99, 10).Value = 3908649.5
$ws.Cells.Item(99, 11).Value = 11112989
$ws.Cells.Item(99, 12).Value = 3908649.5
$ws.Cells.Item(99, 13).Value = -11111491
$ws.Cells.Item(99, 14).Value = -3911645.5
# Row 107: White Oak Lumber
$ws.Cells.Item(107, 8).Value = 1110.9231
$ws.Cells.Item(107, 9).Value = 767.55554
$ws.Cells.Item(107, 10).Value = 1883.5
$ws.Cells.Item(107, 11).Value = 767.55554
$ws.Cells.Item(107, 12).Value = 1883.5
$ws.Cells.Item(107, 13).Value = 1152.44446
$ws.Cells.Item(107, 14).Value = -5723.5
# Row 126: Red Pine Lumber
$ws.Cells.Item(126, 8).Value = 7911060.5
$ws.Cells.Item(126, 9).Value = 11112989
$ws.Cells.Item(126, 10).Value = 3908649.5
$ws.Cells.Item(126, 11).Value = 33338967
$ws.Cells.Item(126, 12).Value = 11725948.5
$ws.Cells.Item(126, 13).Value = -33336497
$ws.Cells.Item(126, 14).Value = -11730888.5

$ws = $wb.Worksheets.Item("CUL")
# Row 25: Apple Tart
$ws.Cells.Item(25, 8).Value = 400
$ws.Cells.Item(25, 9).Value = 400
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 1200
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = -1031
$ws.Cells.Item(25, 14).ClearContents()  # N25
# Row 26: Grape Juice
$ws.Cells.Item(26, 8).Value = 13387.875
$ws.Cells.Item(26, 9).Value = 673.3333
$ws.Cells.Item(26, 10).Value = 21016.6
$ws.Cells.Item(26, 11).Value = 2019.9999
$ws.Cells.Item(26, 12).Value = 63049.8
$ws.Cells.Item(26, 13).Value = -1731.9999
$ws.Cells.Item(26, 14).Value = -63625.8
# Row 30: Apple Tart
$ws.Cells.Item(30, 8).Value = 400
$ws.Cells.Item(30, 9).Value = 400
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 1200
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = -1098
$ws.Cells.Item(30, 14).ClearContents()  # N30
# Row 31: Shepherd's Pie
$ws.Cells.Item(31, 8).Value = 990
$ws.Cells.Item(31, 9).Value = 990
$ws.Cells.Item(31, 11).Value = 2970
$ws.Cells.Item(31, 13).Value = -2682
# Row 50: Rolanberry Cheese
$ws.Cells.Item(50, 8).Value = 1207
$ws.Cells.Item(50, 9).Value = 908.1667
$ws.Cells.Item(50, 10).Value = 3000
$ws.Cells.Item(50, 11).Value = 2724.5001
$ws.Cells.Item(50, 12).Value = 9000
$ws.Cells.Item(50, 13).Value = -2243.5001
$ws.Cells.Item(50, 14).Value = -9962
# Row 53: Rolanberry Cheese
$ws.Cells.Item(53, 8).Value = 1207
$ws.Cells.Item(53, 9).Value = 908.1667
$ws.Cells.Item(53, 10).Value = 3000
$ws.Cells.Item(53, 11).Value = 2724.5001
$ws.Cells.Item(53, 12).Value = 9000
$ws.Cells.Item(53, 13).Value = -2243.5001
$ws.Cells.Item(53, 14).Value = -9962
# Row 81: Frozen Spirits
$ws.Cells.Item(81, 8).Value = 5481.8706
$ws.Cells.Item(81, 10).Value = 7920.6284
$ws.Cells.Item(81, 12).Value = 23761.8852
$ws.Cells.Item(81, 14).Value = -26007.8852
# Row 84: Frozen Spirits
$ws.Cells.Item(84, 8).Value = 5481.8706
$ws.Cells.Item(84, 10).Value = 7920.6284
$ws.Cells.Item(84, 12).Value = 71285.6556
$ws.Cells.Item(84, 14).Value = -82517.6556
# Row 109: Purple Carrot Juice
$ws.Cells.Item(109, 8).Value = 1868.2609
$ws.Cells.Item(109, 9).Value = 935.625
$ws.Cells.Item(109, 10).Value = 4000
$ws.Cells.Item(109, 11).Value = 2806.875
$ws.Cells.Item(109, 12).Value = 12000
$ws.Cells.Item(109, 13).Value = -1766.875
$ws.Cells.Item(109, 14).Value = -14080
# Row 139: Wild Banana Blend
$ws.Cells.Item(139, 8).Value = 4500.75
$ws.Cells.Item(139, 9).Value = 2876.25
$ws.Cells.Item(139, 10).Value = 7749.75
$ws.Cells.Item(139, 11).Value = 8628.75
$ws.Cells.Item(139, 12).Value = 23249.25
$ws.Cells.Item(139, 13).Value = -3488.75
$ws.Cells.Item(139, 14).Value = -33529.25

$ws = $wb.Worksheets.Item("GSM")
# Row 107: Hard Mudstone Whetstone
$ws.Cells.Item(107, 8).Value = 1103.6
$ws.Cells.Item(107, 9).Value = 1016.6667
$ws.Cells.Item(107, 10).Value = 1234
$ws.Cells.Item(107, 11).Value = 1016.6667
$ws.Cells.Item(107, 12).Value = 1234
$ws.Cells.Item(107, 13).Value = 903.3333
$ws.Cells.Item(107, 14).Value = -5074

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Leather
$ws.Cells.Item(7, 8).Value = 3194.28
$ws.Cells.Item(7, 9).Value = 2207.2307
$ws.Cells.Item(7, 11).Value = 2207.2307
$ws.Cells.Item(7, 13).Value = -2095.2307
# Row 82: Dragon Leather
$ws.Cells.Item(82, 8).Value = 3265
$ws.Cells.Item(82, 9).Value = 2775
$ws.Cells.Item(82, 11).Value = 2775
$ws.Cells.Item(82, 13).Value = -2414
# Row 85: Dragon Leather
$ws.Cells.Item(85, 8).Value = 3265
$ws.Cells.Item(85, 9).Value = 2775
$ws.Cells.Item(85, 11).Value = 2775
$ws.Cells.Item(85, 13).Value = -1527
# Row 126: Saiga Leather
$ws.Cells.Item(126, 8).Value = 3194.28
$ws.Cells.Item(126, 9).Value = 2207.2307
$ws.Cells.Item(126, 11).Value = 6621.6921
$ws.Cells.Item(126, 13).Value = -4151.6921
# Row 132: Silver Lobo Leather
$ws.Cells.Item(132, 8).Value = 2048.8333
$ws.Cells.Item(132, 9).Value = 2007.8182
$ws.Cells.Item(132, 11).Value = 6023.4546
$ws.Cells.Item(132, 13).Value = -3493.4546
# Row 136: Br'aax Leather
$ws.Cells.Item(136, 8).Value = 4357.6
$ws.Cells.Item(136, 9).Value = 7245.5
$ws.Cells.Item(136, 11).Value = 21736.5
$ws.Cells.Item(136, 13).Value = -19186.5

$ws = $wb.Worksheets.Item("WVR")
# Row 126: Snow Linen
$ws.Cells.Item(126, 8).Value = 3026.348
$ws.Cells.Item(126, 9).Value = 2420.0667
$ws.Cells.Item(126, 11).Value = 7260.2001
$ws.Cells.Item(126, 13).Value = -4790.2001
